# Update the "想去人数" (want-to-go count) values in the 展览 and 全部类型
# sheets to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6538
$ws1.Range("F15").Value = 3234
$ws1.Range("F17").Value = 203
$ws1.Range("F18").Value = 1888

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6538
$ws4.Range("F16").Value = 3234
$ws4.Range("F18").Value = 203
$ws4.Range("F19").Value = 1888
